$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 214.86667
$ws.Cells.Item(12, 9).Value = 172.58333
$ws.Cells.Item(12, 10).Value = 384
$ws.Cells.Item(12, 11).Value = 172.58333
$ws.Cells.Item(12, 12).Value = 384
$ws.Cells.Item(12, 13).Value = -2.583329999999989
$ws.Cells.Item(12, 14).Value = -724
$ws.Cells.Item(32, 8).Value = 3341.111
$ws.Cells.Item(32, 9).Value = 2759.2
$ws.Cells.Item(32, 10).Value = 4068.5
$ws.Cells.Item(32, 11).Value = 2759.2
$ws.Cells.Item(32, 12).Value = 4068.5
$ws.Cells.Item(32, 13).Value = -2433.2
$ws.Cells.Item(32, 14).Value = -4720.5
$ws.Cells.Item(33, 8).Value = 260.54544
$ws.Cells.Item(33, 9).Value = 181.6
$ws.Cells.Item(33, 10).Value = 1050
$ws.Cells.Item(33, 11).Value = 181.6
$ws.Cells.Item(33, 12).Value = 1050
$ws.Cells.Item(33, 13).Value = 47.40000000000001
$ws.Cells.Item(33, 14).Value = -1508
$ws.Cells.Item(40, 8).Value = 2723043.2
$ws.Cells.Item(40, 9).Value = 7999.222
$ws.Cells.Item(40, 10).Value = 4468428.5
$ws.Cells.Item(40, 11).Value = 7999.222
$ws.Cells.Item(40, 12).Value = 4468428.5
$ws.Cells.Item(40, 13).Value = -7824.222
$ws.Cells.Item(40, 14).Value = -4468778.5
$ws.Cells.Item(58, 8).Value = 1127.4445
$ws.Cells.Item(58, 9).Value = 661.4286
$ws.Cells.Item(58, 10).Value = 2758.5
$ws.Cells.Item(58, 11).Value = 1984.2858
$ws.Cells.Item(58, 12).Value = 8275.5
$ws.Cells.Item(58, 13).Value = -1834.2858
$ws.Cells.Item(58, 14).Value = -8575.5
$ws.Cells.Item(61, 8).Value = 1605.75
$ws.Cells.Item(61, 9).Value = 1605.75
$ws.Cells.Item(61, 11).Value = 4817.25
$ws.Cells.Item(61, 13).Value = -4645.25
$ws.Cells.Item(70, 8).Value = 3779.7334
$ws.Cells.Item(70, 9).Value = 1660
$ws.Cells.Item(70, 10).Value = 4839.6
$ws.Cells.Item(70, 11).Value = 4980
$ws.Cells.Item(70, 12).Value = 14518.8
$ws.Cells.Item(70, 13).Value = -4710
$ws.Cells.Item(70, 14).Value = -15058.8
$ws.Cells.Item(73, 8).Value = 3779.7334
$ws.Cells.Item(73, 9).Value = 1660
$ws.Cells.Item(73, 10).Value = 4839.6
$ws.Cells.Item(73, 11).Value = 4980
$ws.Cells.Item(73, 12).Value = 14518.8
$ws.Cells.Item(73, 13).Value = -4044
$ws.Cells.Item(73, 14).Value = -16390.8
$ws.Cells.Item(112, 8).Value = 10041.077
$ws.Cells.Item(112, 9).Value = 697.2
$ws.Cells.Item(112, 11).Value = 2091.6
$ws.Cells.Item(112, 13).Value = -983.6000000000004
$ws.Cells.Item(125, 8).Value = 3036.25
$ws.Cells.Item(125, 9).Value = 1044.6666
$ws.Cells.Item(125, 10).Value = 3700.111
$ws.Cells.Item(125, 11).Value = 9401.999400000001
$ws.Cells.Item(125, 12).Value = 33300.999
$ws.Cells.Item(125, 13).Value = -6941.999400000001
$ws.Cells.Item(125, 14).Value = -38220.999
$ws.Cells.Item(129, 8).Value = 1730.9131
$ws.Cells.Item(129, 9).Value = 1834.25
$ws.Cells.Item(129, 10).Value = 1494.7142
$ws.Cells.Item(129, 11).Value = 5502.75
$ws.Cells.Item(129, 12).Value = 4484.142599999999
$ws.Cells.Item(129, 13).Value = -502.75
$ws.Cells.Item(129, 14).Value = -14484.1426
$ws.Cells.Item(137, 8).Value = 984732.0600000001
$ws.Cells.Item(137, 9).Value = 965.9286
$ws.Cells.Item(137, 10).Value = 1709612.4
$ws.Cells.Item(137, 11).Value = 2897.7858
$ws.Cells.Item(137, 12).Value = 5128837.199999999
$ws.Cells.Item(137, 13).Value = -347.7857999999997
$ws.Cells.Item(137, 14).Value = -5133937.199999999
$ws.Cells.Item(141, 8).Value = 2650
$ws.Cells.Item(141, 9).Value = 2195
$ws.Cells.Item(141, 11).Value = 6585
$ws.Cells.Item(141, 13).Value = -1405

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(21, 8).Value = 1598
$ws.Cells.Item(21, 9).Value = 248.75
$ws.Cells.Item(21, 10).Value = 6995
$ws.Cells.Item(21, 11).Value = 248.75
$ws.Cells.Item(21, 12).Value = 6995
$ws.Cells.Item(21, 13).Value = 125.25
$ws.Cells.Item(21, 14).Value = -7743
$ws.Cells.Item(30, 8).Value = 2351.6
$ws.Cells.Item(30, 9).Value = 1769.3334
$ws.Cells.Item(30, 10).Value = 3225
$ws.Cells.Item(30, 11).Value = 1769.3334
$ws.Cells.Item(30, 12).Value = 3225
$ws.Cells.Item(30, 13).Value = -1619.3334
$ws.Cells.Item(30, 14).Value = -3525
$ws.Cells.Item(32, 8).Value = 6292692
$ws.Cells.Item(32, 9).Value = 7249965
$ws.Cells.Item(32, 10).Value = 2043.4286
$ws.Cells.Item(32, 11).Value = 7249965
$ws.Cells.Item(32, 12).Value = 2043.4286
$ws.Cells.Item(32, 13).Value = -7249678
$ws.Cells.Item(32, 14).Value = -2617.4286
$ws.Cells.Item(35, 8).Value = 2799.5
$ws.Cells.Item(35, 9).Value = 2799.5
$ws.Cells.Item(35, 11).Value = 2799.5
$ws.Cells.Item(35, 13).Value = -2393.5
$ws.Cells.Item(102, 8).Value = 37463.3
$ws.Cells.Item(102, 9).Value = 37463.3
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 37463.3
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -35841.3
$ws.Cells.Item(102, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1178.7858
$ws.Cells.Item(20, 9).Value = 964.4286
$ws.Cells.Item(20, 10).Value = 1393.1428
$ws.Cells.Item(20, 11).Value = 964.4286
$ws.Cells.Item(20, 12).Value = 1393.1428
$ws.Cells.Item(20, 13).Value = -717.4286
$ws.Cells.Item(20, 14).Value = -1887.1428
$ws.Cells.Item(94, 8).Value = 801.6667
$ws.Cells.Item(94, 9).Value = 739.25
$ws.Cells.Item(94, 11).Value = 739.25
$ws.Cells.Item(94, 13).Value = -288.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3233.3333
$ws.Cells.Item(16, 9).Value = 2850
$ws.Cells.Item(16, 11).Value = 2850
$ws.Cells.Item(16, 13).Value = -2563
$ws.Cells.Item(31, 8).Value = 77988.414
$ws.Cells.Item(31, 9).Value = 121331.16
$ws.Cells.Item(31, 10).Value = 21642.85
$ws.Cells.Item(31, 11).Value = 121331.16
$ws.Cells.Item(31, 12).Value = 21642.85
$ws.Cells.Item(31, 13).Value = -121036.16
$ws.Cells.Item(31, 14).Value = -22232.85
$ws.Cells.Item(34, 8).Value = 77988.414
$ws.Cells.Item(34, 9).Value = 121331.16
$ws.Cells.Item(34, 10).Value = 21642.85
$ws.Cells.Item(34, 11).Value = 121331.16
$ws.Cells.Item(34, 12).Value = 21642.85
$ws.Cells.Item(34, 13).Value = -121129.16
$ws.Cells.Item(34, 14).Value = -22046.85
$ws.Cells.Item(58, 8).Value = 688995.9399999999
$ws.Cells.Item(58, 9).Value = 1123593.6
$ws.Cells.Item(58, 10).Value = 6056.7144
$ws.Cells.Item(58, 11).Value = 1123593.6
$ws.Cells.Item(58, 12).Value = 6056.7144
$ws.Cells.Item(58, 13).Value = -1123390.6
$ws.Cells.Item(58, 14).Value = -6462.7144
$ws.Cells.Item(107, 8).Value = 1032.5714
$ws.Cells.Item(107, 9).Value = 961.63635
$ws.Cells.Item(107, 11).Value = 961.63635
$ws.Cells.Item(107, 13).Value = 958.36365
$ws.Cells.Item(113, 8).Value = 3233.3333
$ws.Cells.Item(113, 9).Value = 2850
$ws.Cells.Item(113, 11).Value = 2850
$ws.Cells.Item(113, 13).Value = -680
$ws.Cells.Item(134, 8).Value = 7535091
$ws.Cells.Item(134, 9).Value = 56776.332
$ws.Cells.Item(134, 10).Value = 18752562
$ws.Cells.Item(134, 11).Value = 170328.996
$ws.Cells.Item(134, 12).Value = 56257686
$ws.Cells.Item(134, 13).Value = -167793.996
$ws.Cells.Item(134, 14).Value = -56262756
$ws.Cells.Item(136, 8).Value = 688995.9399999999
$ws.Cells.Item(136, 9).Value = 1123593.6
$ws.Cells.Item(136, 10).Value = 6056.7144
$ws.Cells.Item(136, 11).Value = 3370780.8
$ws.Cells.Item(136, 12).Value = 18170.1432
$ws.Cells.Item(136, 13).Value = -3368230.8
$ws.Cells.Item(136, 14).Value = -23270.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 14288018
$ws.Cells.Item(121, 9).Value = 33333536
$ws.Cells.Item(121, 10).Value = 3880
$ws.Cells.Item(121, 11).Value = 100000608
$ws.Cells.Item(121, 12).Value = 11640
$ws.Cells.Item(121, 13).Value = -99999298
$ws.Cells.Item(121, 14).Value = -14260

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(38, 8).Value = 24021.5
$ws.Cells.Item(38, 10).Value = 24021.5
$ws.Cells.Item(38, 12).Value = 24021.5
$ws.Cells.Item(38, 14).Value = -24947.5
$ws.Cells.Item(46, 8).Value = 8593.6
$ws.Cells.Item(46, 10).Value = 40046
$ws.Cells.Item(46, 12).Value = 40046
$ws.Cells.Item(46, 14).Value = -40358
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 246198.23
$ws.Cells.Item(80, 9).Value = 368019.16
$ws.Cells.Item(80, 11).Value = 368019.16
$ws.Cells.Item(80, 13).Value = -367021.16
$ws.Cells.Item(83, 8).Value = 246198.23
$ws.Cells.Item(83, 9).Value = 368019.16
$ws.Cells.Item(83, 11).Value = 1840095.8
$ws.Cells.Item(83, 13).Value = -1835103.8
$ws.Cells.Item(132, 8).Value = 18078814
$ws.Cells.Item(132, 9).Value = 25956140
$ws.Cells.Item(132, 10).Value = 7305.2354
$ws.Cells.Item(132, 11).Value = 77868420
$ws.Cells.Item(132, 12).Value = 21915.7062
$ws.Cells.Item(132, 13).Value = -77865890
$ws.Cells.Item(132, 14).Value = -26975.7062

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3960.1738
$ws.Cells.Item(7, 9).Value = 3749.111
$ws.Cells.Item(7, 11).Value = 3749.111
$ws.Cells.Item(7, 13).Value = -3637.111
$ws.Cells.Item(40, 8).Value = 3538.0688
$ws.Cells.Item(40, 9).Value = 3382.8262
$ws.Cells.Item(40, 10).Value = 4133.1665
$ws.Cells.Item(40, 11).Value = 3382.8262
$ws.Cells.Item(40, 12).Value = 4133.1665
$ws.Cells.Item(40, 13).Value = -3246.8262
$ws.Cells.Item(40, 14).Value = -4405.1665
$ws.Cells.Item(126, 8).Value = 3960.1738
$ws.Cells.Item(126, 9).Value = 3749.111
$ws.Cells.Item(126, 11).Value = 11247.333
$ws.Cells.Item(126, 13).Value = -8777.332999999999
$ws.Cells.Item(132, 8).Value = 895074.7
$ws.Cells.Item(132, 9).Value = 1201887.9
$ws.Cells.Item(132, 10).Value = 5316.7
$ws.Cells.Item(132, 11).Value = 3605663.7
$ws.Cells.Item(132, 12).Value = 15950.1
$ws.Cells.Item(132, 13).Value = -3603133.7
$ws.Cells.Item(132, 14).Value = -21010.1
$ws.Cells.Item(136, 8).Value = 57848.39
$ws.Cells.Item(136, 9).Value = 2308.7273
$ws.Cells.Item(136, 11).Value = 6926.1819
$ws.Cells.Item(136, 13).Value = -4376.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2554.0667
$ws.Cells.Item(122, 9).Value = 2017.3334
$ws.Cells.Item(122, 11).Value = 6052.0002
$ws.Cells.Item(122, 13).Value = -3602.0002
$ws.Cells.Item(136, 8).Value = 1451936.4
$ws.Cells.Item(136, 9).Value = 1643194.9
$ws.Cells.Item(136, 11).Value = 4929584.699999999
$ws.Cells.Item(136, 13).Value = -4927034.699999999
